# Fix the weekly schedule grid on Sheet1 of MCT-2A so each turn lasts the
# correct 6 hours: shift the afternoon break ("Almoço") and every slot after
# it down by one 50-minute period, and correct the morning/afternoon class
# assignments that had drifted out of place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "MCT-2A"
$ws.Range("B1").Value = "segunda"
$ws.Range("C1").Value = "terça"
$ws.Range("D1").Value = "quarta"
$ws.Range("E1").Value = "quinta"
$ws.Range("F1").Value = "sexta"

$ws.Range("A2").Value = "7:00"
$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = "Josivaldo Ferreira-Circuitos Elétricos 2"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "-"
$ws.Range("F2").Value = "-"

$ws.Range("A3").Value = "7:50"
$ws.Range("B3").Value = "Lucas Ferreira-Sistemas digitais"
$ws.Range("C3").Value = "Josivaldo Ferreira-Circuitos Elétricos 2"
$ws.Range("D3").Value = "Pedro Francisco-MTRM"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "-"

$ws.Range("A4").Value = "8:40"
$ws.Range("B4").Value = "Lucas Ferreira-Sistemas digitais"
$ws.Range("C4").Value = "Andre Lucca-Acionamentos"
$ws.Range("D4").Value = "Pedro Francisco-MTRM"
$ws.Range("E4").Value = "Andre Barros-EAP"
$ws.Range("F4").Value = "-"

$ws.Range("A5").Value = "9:30"
$ws.Range("B5").Value = "Intervalo"
$ws.Range("C5").Value = "Intervalo"
$ws.Range("D5").Value = "Intervalo"
$ws.Range("E5").Value = "Intervalo"
$ws.Range("F5").Value = "Intervalo"

$ws.Range("A6").Value = "9:50"
$ws.Range("B6").Value = "-"
$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "Josivaldo Ferreira-Programação"
$ws.Range("F6").Value = "João Rodrigues-CAD"

$ws.Range("A7").Value = "10:40"
$ws.Range("B7").Value = "Andre Barros-EAP"
$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "-"
$ws.Range("F7").Value = "João Rodrigues-CAD"

$ws.Range("A8").Value = "11:30"
$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "-"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "-"
$ws.Range("F8").Value = "-"

$ws.Range("A9").Value = "12:20"
$ws.Range("B9").Value = "Almoço"
$ws.Range("C9").Value = "Almoço"
$ws.Range("D9").Value = "Almoço"
$ws.Range("E9").Value = "Almoço"
$ws.Range("F9").Value = "Almoço"

$ws.Range("A10").Value = "13:00"
$ws.Range("B10").Value = "-"
$ws.Range("C10").Value = "-"
$ws.Range("D10").Value = "-"
$ws.Range("E10").Value = "-"
$ws.Range("F10").Value = "-"

$ws.Range("A11").Value = "13:50"
$ws.Range("B11").Value = "-"
$ws.Range("C11").Value = "-"
$ws.Range("D11").Value = "-"
$ws.Range("E11").Value = "-"
$ws.Range("F11").Value = "-"

$ws.Range("A12").Value = "14:40"
$ws.Range("B12").Value = "-"
$ws.Range("C12").Value = "-"
$ws.Range("D12").Value = "-"
$ws.Range("E12").Value = "-"
$ws.Range("F12").Value = "-"

$ws.Range("A13").Value = "15:30"
$ws.Range("B13").Value = "Intervalo"
$ws.Range("C13").Value = "Intervalo"
$ws.Range("D13").Value = "Intervalo"
$ws.Range("E13").Value = "Intervalo"
$ws.Range("F13").Value = "Intervalo"

$ws.Range("A14").Value = "15:50"
$ws.Range("B14").Value = "-"
$ws.Range("C14").Value = "-"
$ws.Range("D14").Value = "-"
$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "-"

$ws.Range("A15").Value = "16:40"
$ws.Range("B15").Value = "-"
$ws.Range("C15").Value = "-"
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"

$ws.Range("A16").Value = "17:30"
$ws.Range("B16").Value = "-"
$ws.Range("C16").Value = "-"
$ws.Range("D16").Value = "-"
$ws.Range("E16").Value = "-"
$ws.Range("F16").Value = "-"

$ws.Range("A17").Value = "18:20"
$ws.Range("B17:F17").Value = ""
# Touch formatting so the (intentionally blank) cells are still materialized
# in the saved sheet, matching the target where B17:F17 exist as empty cells.
$ws.Range("B17:F17").Font.Bold = $false
